$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(-19.24488640892737, 1.951795554027346, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737)
    3 = @(-19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, 1.210164828773727, -19.24488640892737, -19.24488640892737)
    4 = @(-19.24488640892737, 1.991100329374781, 1.662929202594206, -19.24488640892737, 3.392532037942584, -19.24488640892737, -19.24488640892737, -19.24488640892737, 0.972712813498666, -19.24488640892737)
    5 = @(-19.24488640892737, 1.692542730910142, -19.24488640892737, -19.24488640892737, -19.24488640892737, 2.815193665098422, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737)
    6 = @(-19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737)
    7 = @(2.454464486391756, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737)
    8 = @(-19.24488640892737, -19.24488640892737, -19.24488640892737, 1.779783224183466, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737)
    9 = @(3.859857736223567, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737)
    10 = @(-19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, 1.845097396062276, -19.24488640892737, 2.186291330425833)
    11 = @(-19.24488640892737, -19.24488640892737, -19.24488640892737, 2.922562670718259, -19.24488640892737, 2.864504619980575, -19.24488640892737, -19.24488640892737, -19.24488640892737, 1.901963767133809)
    12 = @(-19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737)
    13 = @(-19.24488640892737, -19.24488640892737, -19.24488640892737, 2.537289984738388, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, 1.684011425394473, 1.795508338880507)
    14 = @(-19.24488640892737, -19.24488640892737, 1.515720661165363, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, 2.070953605611128)
    15 = @(-19.24488640892737, -19.24488640892737, 1.771781119697573, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737)
    16 = @(-19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, 1.771100930737816, -19.24488640892737)
    17 = @(-19.24488640892737, 2.128888399977065, 1.81269091617183, -19.24488640892737, -19.24488640892737, -19.24488640892737, 4.321925888862493, 1.980488742300141, 2.591332339305997, -19.24488640892737)
    18 = @(-19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, 2.017367328457969, 2.428826811589039, -19.24488640892737)
    19 = @(-19.24488640892737, -19.24488640892737, 2.050100535469548, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737, 1.767067630065589, -19.24488640892737, -19.24488640892737)
    20 = @(-19.24488640892737, 1.026573199356744, 1.538489957650775, -19.24488640892737, 3.247685869505069, -19.24488640892737, -19.24488640892737, 1.430413398008828, -19.24488640892737, 2.013783705950754)
    21 = @(-19.24488640892737, 1.329283005729398, -19.24488640892737, 1.668725455978895, -19.24488640892737, 2.505631269373283, -19.24488640892737, -19.24488640892737, -19.24488640892737, -19.24488640892737)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, $i + 2).Value = $vals[$i]
    }
}
